# Generate Report for Handoff
# Replaces the old GUID-based file name "b89d49ee-48f5-4d6b-a294-caf147eaec38"
# with the newly generated one "35c6a7cd-6bef-496f-888f-e35934852f7c" across
# all three sheets, updates the hyperlink display text (keeping the original
# hyperlink target), refreshes the xlf hash-qualified target file names, and
# bumps the handoff/handback timestamps to reflect the new report run.

$wb = $excel.ActiveWorkbook

$oldGuid = "b89d49ee-48f5-4d6b-a294-caf147eaec38"
$newGuid = "35c6a7cd-6bef-496f-888f-e35934852f7c"

$oldHash = "486d7c5e36487a7dd6db2df4cca62d4733155557"
$newHash = "33998104ed7215556f29712780f9becdee53c75f"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e9da94a4204f2d2a78a36d6086d81b199138c546/e2e/$oldGuid.md"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newGuid.md")
$wsOverview.Range("G2").Value = "2016-08-19 04:52:01"

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-19 04:51:55"

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-19 04:52:01"
